# Apply the "train on assignment 2 and test on assignment 1" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column B: fill in plagiarism-check results for rows 2-20 (assignment2-1 .. assignment2-19) ---
$bValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
}

foreach ($row in $bValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $bValues[$row]
}

# --- Column A: clear out the old assignment2-20 .. assignment2-54 file names (rows 21-55) ---
$ws.Range("A21:A55").ClearContents()

# --- Update the active selection to reflect where the author ended up clicking ---
$ws.Range("H16").Select()
